# Update code for report co so
# Insert 3 new columns (ED:EF) for a new "Ngày 7" select property group,
# shifting the existing "Ngày 1" and "Ngày 4" select columns three places
# to the right (ED->EG, EE->EH, EF->EI, EG->EJ, EH->EK, EI->EL), then
# populate the freed ED:EF columns with the "Ngày 7" header/data and bump
# a handful of last_edited_time stamps.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert three whole columns starting at ED, shifting ED:EI (and
#    everything to their right) over to EG:EL. This also grows the
#    sheet dimension from A1:EI14 to A1:EL14 automatically.
$ws.Range("ED1:EF1").EntireColumn.Insert()

# 2. Fill in the header row for the new "Ngày 7" select property.
$ws.Range("ED1").Value = "properties.Ngày 7.select.id"
$ws.Range("EE1").Value = "properties.Ngày 7.select.name"
$ws.Range("EF1").Value = "properties.Ngày 7.select.color"

# 3. Fill in the new "Ngày 7" select data for each data row.
$ws.Range("ED6").Value = "VRLp"
$ws.Range("EE6").Value = "Nghỉ có phép"
$ws.Range("EF6").Value = "blue"

$ws.Range("ED7").Value = "DjwF"
$ws.Range("EE7").Value = "Đầy đủ"
$ws.Range("EF7").Value = "pink"

$ws.Range("ED9").Value = "DjwF"
$ws.Range("EE9").Value = "Đầy đủ"
$ws.Range("EF9").Value = "pink"

$ws.Range("ED10").Value = "DjwF"
$ws.Range("EE10").Value = "Đầy đủ"
$ws.Range("EF10").Value = "pink"

$ws.Range("ED14").Value = "DjwF"
$ws.Range("EE14").Value = "Đầy đủ"
$ws.Range("EF14").Value = "pink"

# 4. Bump last_edited_time for the rows that were touched by this update.
$ws.Range("D6").Value = "2024-07-07T12:43:00.000Z"
$ws.Range("D7").Value = "2024-07-07T12:43:00.000Z"
$ws.Range("D9").Value = "2024-07-07T12:43:00.000Z"
$ws.Range("D10").Value = "2024-07-07T12:43:00.000Z"
$ws.Range("D14").Value = "2024-07-07T12:42:00.000Z"
